# Auto-generated edit script to apply cell value changes per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K18").Value = 2133.3333
$ws.Range("H18").Value = 2133.3333
$ws.Range("I18").Value = 2133.3333
$ws.Range("M18").Value = -1849.3333
$ws.Range("K32").Value = 2165.3333
$ws.Range("M32").Value = -1839.3333
$ws.Range("I32").Value = 2165.3333
$ws.Range("H32").Value = 2285
$ws.Range("I94").Value = 10035
$ws.Range("M94").Value = -9584
$ws.Range("H94").Value = 10035
$ws.Range("K94").Value = 10035
$ws.Range("I111").Value = 2238
$ws.Range("N111").Value = -13127
$ws.Range("M111").Value = -3647
$ws.Range("K111").Value = 6714
$ws.Range("L111").Value = 6993
$ws.Range("H111").Value = 2263.3635
$ws.Range("J111").Value = 2331
$ws.Range("J112").Value = 2555.111
$ws.Range("H112").Value = 2499.6
$ws.Range("L112").Value = 7665.333
$ws.Range("N112").Value = -9881.332999999999
$ws.Range("I125").Value = 983.6667
$ws.Range("H125").Value = 2989
$ws.Range("K125").Value = 8853.0003
$ws.Range("M125").Value = -6393.0003
$ws.Range("K132").Value = 6741.3999
$ws.Range("I132").Value = 2247.1333
$ws.Range("H132").Value = 2335.7058
$ws.Range("M132").Value = -4211.3999
$ws.Range("J138").Value = 4258.4653
$ws.Range("H138").Value = 3716.83
$ws.Range("N138").Value = -23055.3959
$ws.Range("M138").Value = 976.6000000000004
$ws.Range("K138").Value = 4163.4
$ws.Range("L138").Value = 12775.3959
$ws.Range("I138").Value = 1387.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M74").Value = -155.375
$ws.Range("K74").Value = 1029.375
$ws.Range("H74").Value = 1318.1923
$ws.Range("I74").Value = 1029.375
$ws.Range("M77").Value = -778.875
$ws.Range("K77").Value = 5146.875
$ws.Range("I77").Value = 1029.375
$ws.Range("H77").Value = 1318.1923
$ws.Range("N109").Value = -50774
$ws.Range("L109").Value = 48000
$ws.Range("J109").Value = 48000
$ws.Range("H109").Value = 48000
$ws.Range("K132").Value = 4642.5
$ws.Range("I132").Value = 1547.5
$ws.Range("H132").Value = 1547.5
$ws.Range("M132").Value = -2112.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I20").Value = 11699.8
$ws.Range("M20").Value = -11452.8
$ws.Range("H20").Value = 11699.8
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("K20").Value = 11699.8
$ws.Range("K99").Value = 2251.3333
$ws.Range("M99").Value = -753.3332999999998
$ws.Range("H99").Value = 2583.7778
$ws.Range("I99").Value = 2251.3333
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -348.1579999999999
$ws.Range("N105").ClearContents()
$ws.Range("I105").Value = 2095.158
$ws.Range("H105").Value = 2095.158
$ws.Range("K105").Value = 2095.158
$ws.Range("I134").Value = 2375.5833
$ws.Range("M134").Value = -4591.749899999999
$ws.Range("K134").Value = 7126.749899999999
$ws.Range("H134").Value = 2667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K16").Value = 2266.0417
$ws.Range("I16").Value = 2266.0417
$ws.Range("M16").Value = -1979.0417
$ws.Range("H16").Value = 2365.3103
$ws.Range("H22").Value = 500.2
$ws.Range("J22").Value = 634
$ws.Range("N22").Value = -1334
$ws.Range("L22").Value = 634
$ws.Range("I58").Value = 1605.2222
$ws.Range("J58").Value = 5218.231
$ws.Range("K58").Value = 1605.2222
$ws.Range("L58").Value = 5218.231
$ws.Range("M58").Value = -1402.2222
$ws.Range("H58").Value = 3740.182
$ws.Range("N58").Value = -5624.231
$ws.Range("J112").Value = 0
$ws.Range("H112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("I113").Value = 2266.0417
$ws.Range("M113").Value = -96.04170000000022
$ws.Range("K113").Value = 2266.0417
$ws.Range("H113").Value = 2365.3103
$ws.Range("K122").Value = 8057.0625
$ws.Range("I122").Value = 2685.6875
$ws.Range("M122").Value = -5607.0625
$ws.Range("H122").Value = 2824.9
$ws.Range("I134").Value = 1761.3846
$ws.Range("N134").Value = -14664.4614
$ws.Range("M134").Value = -2749.1538
$ws.Range("L134").Value = 9594.4614
$ws.Range("K134").Value = 5284.1538
$ws.Range("H134").Value = 2240.3076
$ws.Range("J134").Value = 3198.1538
$ws.Range("K136").Value = 4815.6666
$ws.Range("L136").Value = 15654.693
$ws.Range("M136").Value = -2265.6666
$ws.Range("I136").Value = 1605.2222
$ws.Range("H136").Value = 3740.182
$ws.Range("J136").Value = 5218.231
$ws.Range("N136").Value = -20754.693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J43").Value = 475
$ws.Range("N43").Value = -1653
$ws.Range("L43").Value = 1425
$ws.Range("H43").Value = 475
$ws.Range("H68").Value = 1186.1666
$ws.Range("K68").Value = 3343.6665
$ws.Range("M68").Value = -2532.6665
$ws.Range("J68").Value = 1401
$ws.Range("I68").Value = 1114.5555
$ws.Range("L68").Value = 4203
$ws.Range("N68").Value = -5825
$ws.Range("J71").Value = 1401
$ws.Range("M71").Value = -5974.9995
$ws.Range("I71").Value = 1114.5555
$ws.Range("K71").Value = 10030.9995
$ws.Range("L71").Value = 12609
$ws.Range("H71").Value = 1186.1666
$ws.Range("N71").Value = -20721
$ws.Range("L103").Value = 9498.999899999999
$ws.Range("H103").Value = 3166.3333
$ws.Range("N103").Value = -11256.9999
$ws.Range("J103").Value = 3166.3333
$ws.Range("M137").Value = -3007.7145
$ws.Range("K137").Value = 8107.7145
$ws.Range("J137").Value = 7370.8
$ws.Range("N137").Value = -32312.4
$ws.Range("L137").Value = 22112.4
$ws.Range("I137").Value = 2702.5715
$ws.Range("H137").Value = 3931.0527

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I126").Value = 3518.6
$ws.Range("H126").Value = 4260.7
$ws.Range("M126").Value = -8085.799999999999
$ws.Range("K126").Value = 10555.8
$ws.Range("K132").Value = 6025.9998
$ws.Range("I132").Value = 2008.6666
$ws.Range("H132").Value = 5133.3335
$ws.Range("M132").Value = -3495.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I55").Value = 183.52942
$ws.Range("H55").Value = 195.1579
$ws.Range("K55").Value = 183.52942
$ws.Range("M55").Value = -10.52941999999999
$ws.Range("K61").Value = 4099.5
$ws.Range("J61").Value = 2338
$ws.Range("N61").Value = -2742
$ws.Range("I61").Value = 4099.5
$ws.Range("H61").Value = 3512.3333
$ws.Range("M61").Value = -3897.5
$ws.Range("L61").Value = 2338
$ws.Range("J82").Value = 2208.7144
$ws.Range("N82").Value = -2930.7144
$ws.Range("L82").Value = 2208.7144
$ws.Range("H82").Value = 3244.238
$ws.Range("J85").Value = 2208.7144
$ws.Range("L85").Value = 2208.7144
$ws.Range("N85").Value = -4704.7144
$ws.Range("H85").Value = 3244.238
$ws.Range("M93").Value = -65.25
$ws.Range("I93").Value = 1313.25
$ws.Range("K93").Value = 1313.25
$ws.Range("H93").Value = 1313.25
$ws.Range("I113").Value = 4099.5
$ws.Range("M113").Value = -1929.5
$ws.Range("K113").Value = 4099.5
$ws.Range("L113").Value = 2338
$ws.Range("H113").Value = 3512.3333
$ws.Range("N113").Value = -6678
$ws.Range("J113").Value = 2338
$ws.Range("K122").Value = 74997
$ws.Range("I122").Value = 24999
$ws.Range("M122").Value = -72547
$ws.Range("H122").Value = 16499.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1883.4546
$ws.Range("L81").Value = 1400
$ws.Range("K81").Value = 4003.6
$ws.Range("J81").Value = 700
$ws.Range("I81").Value = 2001.8
$ws.Range("M81").Value = -2942.6
$ws.Range("N81").Value = -3522
$ws.Range("M84").Value = -14714
$ws.Range("N84").Value = -17608
$ws.Range("J84").Value = 700
$ws.Range("H84").Value = 1883.4546
$ws.Range("L84").Value = 7000
$ws.Range("I84").Value = 2001.8
$ws.Range("K84").Value = 20018
$ws.Range("K96").Value = 1420.5
$ws.Range("J96").Value = 1429
$ws.Range("M96").Value = -47.5
$ws.Range("I96").Value = 1420.5
$ws.Range("N96").Value = -4175
$ws.Range("L96").Value = 1429
$ws.Range("H96").Value = 1425.2222
$ws.Range("K100").Value = 2873.7778
$ws.Range("N100").Value = -5376
$ws.Range("I100").Value = 1436.8889
$ws.Range("J100").Value = 2147
$ws.Range("L100").Value = 4294
$ws.Range("H100").Value = 1614.4166
$ws.Range("M100").Value = -2332.7778
$ws.Range("K122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H122").Value = 3500
$ws.Range("K132").Value = 4162.3638
$ws.Range("I132").Value = 1387.4546
$ws.Range("L132").Value = 4653.75
$ws.Range("J132").Value = 1551.25
$ws.Range("H132").Value = 1431.1333
$ws.Range("M132").Value = -1632.3638
$ws.Range("N132").Value = -9713.75
$ws.Range("K136").Value = 5141.6001
$ws.Range("L136").Value = 14997
$ws.Range("M136").Value = -2591.6001
$ws.Range("I136").Value = 1713.8667
$ws.Range("H136").Value = 2012.5151
$ws.Range("J136").Value = 4999
$ws.Range("N136").Value = -20097
$ws.Range("N139").Value = -103495
$ws.Range("H139").Value = 48061.43
$ws.Range("L139").Value = 93215
$ws.Range("J139").Value = 93215
